# Daylight Savings Time update for the Shadow Lookup Table.
#
# The event names in column A are renumbered/relabeled (DST shifted every
# event's clock time, so the old "<Event> <time>" labels no longer apply -
# they are replaced with simple sequential labels). Two rows (the old
# "Raid the Vault 7PM" / "Shadow Assembly 7PM" pair) also swap position
# relative to each other in terms of which shared-string slot they use,
# which nets out to the "Raid the Vault" / "Shadow Assembly" labels
# trading places between rows 3 and 4.
#
# Final column A text (row by row) after the edit:
#   A1  Event            (unchanged header)
#   A2  Raid the Vault 1
#   A3  Raid the Vault 2
#   A4  Shadow Assembly
#   A5  Battlegrounds 1
#   A6  Battlegrounds 2
#   A7  Battlegrounds 3
#   A8  Battlegrounds 4
#   A9  Shadow Lottery 1
#   A10 Shadow Lottery 2
#   A11 Shadow Lottery 3
#   A12 Shadow War
#   A13 Rite of Exile

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Raid the Vault 1"
$ws.Range("A3").Value  = "Raid the Vault 2"
$ws.Range("A4").Value  = "Shadow Assembly"
$ws.Range("A5").Value  = "Battlegrounds 1"
$ws.Range("A6").Value  = "Battlegrounds 2"
$ws.Range("A7").Value  = "Battlegrounds 3"
$ws.Range("A8").Value  = "Battlegrounds 4"
$ws.Range("A9").Value  = "Shadow Lottery 1"
$ws.Range("A10").Value = "Shadow Lottery 2"
$ws.Range("A11").Value = "Shadow Lottery 3"
$ws.Range("A12").Value = "Shadow War"
$ws.Range("A13").Value = "Rite of Exile"

# The author's active cell/selection moved from A5 to C8 before saving.
$ws.Range("C8").Select()
